# logboek.xlsx update: add two new logged-work entries (rows 16 & 17)
# and move the active-cell selection to A17, matching the commit:
# "api op geboden, toeveogen van bod weg doen, algemene functie problemen op gelost"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: new log entry - 18 Dec 2024, 3 worked hours (column D)
$ws.Range("A16").Value = "na login pagina error weg gehaald, overlopen van alle funcionaliteiten en de functies die niet meer werken gefixt, rating laten weer geven en bod weg doen voor klusjes man"
$ws.Range("B16").Value = "12/18/2024"
$ws.Range("D16").Value = 3

# Row 17: new log entry - 19 Dec 2024, 1 worked hour (column D)
$ws.Range("A17").Value = "geboden functie via api"
$ws.Range("B17").Value = "12/19/2024"
$ws.Range("D17").Value = 1

# Move selection to A17, where the cursor ended up after the edits
$ws.Range("A17").Select()
